$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price increases
$ws.Range("B9").Value = 4940.0

$ws.Range("B42:B47").Value = 730.0

$ws.Range("B72:B74").Value = 4090.0
$ws.Range("B75").Value = 4960.0

$ws.Range("B82:B86").Value = 730.0

$ws.Range("B93:B94").Value = 2300.0
$ws.Range("B95:B97").Value = 1420.0
